$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6128789782524109
$ws.Range("B1").Value = 2.235052824020386
$ws.Range("C1").Value = 6.205644607543945
$ws.Range("D1").Value = 1.620736241340637
$ws.Range("E1").Value = 1.591665029525757
